$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "AfterShip"

# Clear out the old data rows (2-5) and rebuild the header row with the new
# column order / names.
$ws.Cells.Clear()

$ws.Range("A1").Value = "tracking_number"
$ws.Range("B1").Value = "carrier_slug"
$ws.Range("C1").Value = "status_tag"
$ws.Range("D1").Value = "order_id"
$ws.Range("E1").Value = "last_checkpoint_id"
$ws.Range("F1").Value = "last_checkpoint_time"
$ws.Range("G1").Value = "last_checkpoint_location"
$ws.Range("H1").Value = "updated_at"
$ws.Range("I1").Value = "title"
